$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.699.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.270.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.277.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0974"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.48%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.674.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.808.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.291.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0701"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("E32").Value = "  +5.65%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.901"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.90%  "
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  +3.49%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.40%  "
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "256.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0504"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0913"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.70%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0211"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("B50").Value = "Polygon"
$ws.Range("C50").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.373"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("E51").Value = "  +1.35%  "
